# Qual_5_data.xlsx — "created series data tables"
#
# The MINSET_SERIES rows that had been appended at the bottom of the sheet
# (rows 192:316) were missing their INDICATOR_ID (column A) value. Backfill
# them with 58, matching the rest of the table, then (re)apply the sheet's
# sort/autofilter on column A and move the selection/scroll position to
# where the author left off reviewing the newly filled-in rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Backfill INDICATOR_ID (column A) for the rows that were missing it.
#    Rows 2:191 already carried the value 58; rows 192:316 did not.
$ws.Range("A192:A316").Value = 58

# 2) Re-sort/re-apply the autofilter range over the now fully-populated
#    column A (mirrors Data > Sort performed through the table's filter).
$sort = $ws.AutoFilter().Sort()
$sort.SortFields().Clear()
$sort.SortFields().Add2($ws.Range("A38"), $null, 1, $null, $null)
$sort.SetRange($ws.Range("A1:O316"))
$sort.Header = 1
$sort.Apply()

# 3) Leave the view where the author was working: scrolled down with the
#    frozen header still showing, and A203:A204 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 140
$win.ScrollColumn = 1
[void]$ws.Range("A203:A204").Select()
